# Updates the "cryptos" worksheet with refreshed price/volume data.
# Mirrors the upstream GitHub Actions commit that refreshed coin prices
# and 1h volume percentages, including a ShibaInu/Polygon row swap
# (rows 14 and 15).
#
# Cells whose new text would otherwise be auto-parsed by Excel as a
# number (e.g. "1.000", "0.6304") are entered with a leading apostrophe
# so they remain plain text, matching the original inline-string cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.449.25"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.850.80"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'240.94"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").Value = "'0.6304"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.07704"
$ws.Range("E8").Value = "  +2.38%  "
$ws.Range("D9").Value = "'0.2929"
$ws.Range("E9").Value = "  -0.47%  "
$ws.Range("D10").Value = "'24.73"
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("D11").Value = "'0.07742"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").Value = "1.894.40"
$ws.Range("E12").Value = "  +2.25%  "
$ws.Range("E13").Value = "  +1.29%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.6796"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.00001073"
$ws.Range("E15").Value = "  +4.76%  "
$ws.Range("D16").Value = "'83.72"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("D17").Value = "2.156.84"
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").Value = "'6.211"
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("D19").Value = "29.486.37"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").Value = "'228.57"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").Value = "'1.0000"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "'7.465"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "'157.59"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").Value = "'8.417"
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("E29").Value = "  +5.88%  "
$ws.Range("E30").Value = "  +0.62%  "
$ws.Range("D31").Value = "'0.05678"
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").Value = "'4.038"
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("D36").Value = "'0.7035"
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("D37").Value = "'2.586"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").Value = "'2.781"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("D39").Value = "'0.01791"
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("D40").Value = "1.219.61"
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("D41").Value = "'6.556"
$ws.Range("E41").Value = "  +5.26%  "
$ws.Range("D42").Value = "'0.9069"
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").Value = "'101.73"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "'66.43"
$ws.Range("E45").Value = "  +1.43%  "
$ws.Range("E46").Value = "  +1.40%  "
$ws.Range("E47").Value = "  +0.75%  "
$ws.Range("D48").Value = "'0.4023"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("D49").Value = "'9.011"
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("D50").Value = "'1.685"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("D51").Value = "'0.1142"
$ws.Range("E51").Value = "  +2.10%  "
